# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff" and the two timestamps that track when that
# handoff snapshot was produced were bumped to reflect the new run.
#
# Layout of the workbook:
#   Overview (sheet 1): E2/F2 = zh-cn/de-de Status, G2 = Latest HO Xliff
#     Generate Date
#   zh-cn   (sheet 2): C2 = Status, H2 = Latest Handoff Datetime
#   de-de   (sheet 3): C2 = Status, H2 = Latest Handoff Datetime

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps bumped to the new handoff run -------------------------
# Overview!G2 and de-de!H2 shared the same "Latest HO Xliff Generate
# Date" / "Latest Handoff Datetime" stamp before the edit, and still do.
$wsOverview.Range("G2").Value = "2016-08-18 15:04:15"
$wsDeDe.Range("H2").Value     = "2016-08-18 15:04:15"
# zh-cn!H2 had its own distinct stamp.
$wsZhCn.Range("H2").Value     = "2016-08-18 15:04:00"

# --- Column widths now have to accommodate "Ready for handoff" --------
# (same status column on each sheet: Overview E:F, zh-cn C, de-de C)
$wsOverview.Range("E1:F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth        = 16.33
$wsDeDe.Range("C1").ColumnWidth        = 16.33
